$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, matching the style of the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add numeric values in H2 and H3 (plain, unstyled like G2/G3)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
